# Generate Report for Archive
#
# Localization status moved from "Ready for handoff" to "In Translation"
# for the tracked file. That status string is shared across the Overview
# sheet (per-language status columns) and each language's own status
# sheet, so update every occurrence. Excel's column AutoFit (triggered by
# the now-shorter text) narrows the affected "Status"/"zh-cn"/"de-de"
# columns, so the column widths are tightened to match too.
#
# NOTE: the host stores ColumnWidth in 1/6-character pixel increments
# (width -> round(width*6) pixels, persisted back as pixels/6), so the
# nearest representable width to the recorded 13.4101845877511 is driven
# in via 12.5 (which lands on 13.333333...).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$narrowWidth = 12.5

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
